$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'286.83"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'1.72%"
$ws.Range("E2").Style = "Normal"
$ws.Range("E3").Value = "'4.45%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.063"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'0.55%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.06814"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'5.12%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'7.377"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'2.14%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'1.381"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'-0.19%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.9002"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'-2.98%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.1587"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'1.96%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.06954"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'12.00%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.07639"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'0.96%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.02913"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'1.65%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.08977"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'-0.10%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.001614"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'1.74%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.0006421"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'0.84%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.006495"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'7.13%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'3.459"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'0.55%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'3.459"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'2.24%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'2.231"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'0.01%"
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'0.3219"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'0.58%"
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'2.59%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'3.989"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'-1.80%"
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'0.65%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.04474"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'1.60%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.001201"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'1.52%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.004375"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'-0.32%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'0.0001163"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'-6.98%"
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'0.0001608"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'-0.74%"
$ws.Range("E28").Style = "Normal"
$ws.Range("D40").Value = "'0.04280"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'3.02%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.006774"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'2.13%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.1240"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'1.67%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.002196"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'8.19%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.01159"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'-3.99%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.00005704"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'1.81%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D47").Value = "'0.01299"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'-0.21%"
$ws.Range("E47").Style = "Normal"
